$d = $word.ActiveDocument

# --- 1) Insert a new "Meta description" paragraph right after the H1 title ---
# Build the new paragraph after an existing plain ("Normal"-style) body
# paragraph so it naturally picks up plain formatting (no heading style),
# then move it into place right after the title - this avoids stamping the
# paragraph with an explicit style change.
$bodyPara = $d.Paragraphs.Item(3)
$bodyPara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(4)

$metaStart = $metaPara.Range.Start
$metaFull = $d.Range($metaStart, $metaPara.Range.End - 1)
$metaFull.Text = "Meta description: Read our unbiased review of Diwinity, an online slot game with free spins and bonus game features. Play now for free."

$metaBold = $d.Range($metaStart, $metaStart + 16)
$metaBold.Bold = 1

$metaParaFull = $d.Paragraphs.Item(4).Range
$metaParaFull.Cut()

$titlePara = $d.Paragraphs.Item(1)
$insertAt = $titlePara.Range.End
$pasteRange = $d.Range($insertAt, $insertAt)
$pasteRange.Paste()

# --- 2) Remove the duplicated bold "Play Diwinity..." paragraph near the end ---
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Play Diwinity Free Slot Game - Review & Rating" -and $i -ne 1) {
        $p.Range.Delete()
        break
    }
}

# --- 3) Replace the italic meta-description-like closing paragraph with the image prompt ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$promptText = 'Prompt: Create a cartoon-style feature image for Diwinity that features a happy Maya warrior with glasses. The image should convey excitement and adventure, as well as highlight the theme of ancient gods and mythology. The background should feature elements from the game, such as a fountain and columns. The mage should be vibrant and eye-catching, with bold colors and a dynamic pose for the Maya warrior. It should also include the title of the game - "Diwinity" - and any other relevant text that captures the essence of the game.'
$replaceRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$replaceRange.Text = $promptText
